$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" sheet, with
#    the same layout as the other quarterly fund-holding sheets.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}
$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$fundRows = @(
    @("519629", "银河睿利灵活配置混合A", "2.24", "24.70", "1.07", "0.0240", 9),
    @("519630", "银河睿利灵活配置混合C", "1.94", "24.70", "1.07", "0.0208", 9)
)

for ($r = 0; $r -lt $fundRows.Length; $r++) {
    $rowNum = $r + 2
    $data = $fundRows[$r]

    $idxCell = $q1.Cells.Item($rowNum, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    for ($c = 0; $c -lt 6; $c++) {
        $cell = $q1.Cells.Item($rowNum, $c + 2)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c]
    }
    $q1.Cells.Item($rowNum, 8).Value = $data[6]
}

# ------------------------------------------------------------------
# 2. Add a new row to the "总计" (totals) sheet summarising 2022-Q1,
#    pushing the existing rows down and renumbering the index column.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()

# Re-use the formatting already applied to the index column.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 0.04

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4

Write-Host "2022-Q1 data added"
